$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the values of rows 3..16 up into rows 2..15: copy only the values
# (PasteSpecial values) so each destination row keeps its own
# formatting/style banding, while the text moves up one row. Finally
# remove the now-duplicated last row (16).
for ($r = 3; $r -le 16; $r++) {
    $destRow = $r - 1
    $src = $ws.Range("A" + $r + ":F" + $r)
    $dst = $ws.Range("A" + $destRow + ":F" + $destRow)
    $src.Copy()
    $dst.PasteSpecial(-4163)  # xlPasteValues
}

$excel.CutCopyMode = 0
$ws.Rows(16).Delete()
